$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)

# Row 11
$ws.Cells.Item(11, 8).Value = 97  # H11: 0 -> 97
$ws.Cells.Item(11, 9).Value = 97  # I11: 0 -> 97
$ws.Cells.Item(11, 11).Value = 97  # K11: 0 -> 97
$ws.Cells.Item(11, 13).Value = 43  # M11: None -> 43

# Row 21
$ws.Cells.Item(21, 8).Value = 20000  # H21: 18008.5 -> 20000
$ws.Cells.Item(21, 9).Value = 0  # I21: 16017 -> 0
$ws.Cells.Item(21, 11).Value = 0  # K21: 16017 -> 0
$ws.Cells.Item(21, 13).ClearContents()  # M21: -15549 -> (removed)

# Row 23
$ws.Cells.Item(23, 8).Value = 20000  # H23: 18008.5 -> 20000
$ws.Cells.Item(23, 9).Value = 0  # I23: 16017 -> 0
$ws.Cells.Item(23, 11).Value = 0  # K23: 16017 -> 0
$ws.Cells.Item(23, 13).ClearContents()  # M23: -15783 -> (removed)

# Row 80
$ws.Cells.Item(80, 8).Value = 806.3929000000001  # H80: 791.4828 -> 806.3929000000001
$ws.Cells.Item(80, 9).Value = 789.35297  # I80: 766.2778 -> 789.35297
$ws.Cells.Item(80, 11).Value = 2368.05891  # K80: 2298.8334 -> 2368.05891
$ws.Cells.Item(80, 13).Value = -1370.05891  # M80: -1300.8334 -> -1370.05891

# Row 83
$ws.Cells.Item(83, 8).Value = 806.3929000000001  # H83: 791.4828 -> 806.3929000000001
$ws.Cells.Item(83, 9).Value = 789.35297  # I83: 766.2778 -> 789.35297
$ws.Cells.Item(83, 11).Value = 7104.17673  # K83: 6896.500199999999 -> 7104.17673
$ws.Cells.Item(83, 13).Value = -2112.17673  # M83: -1904.500199999999 -> -2112.17673

# Row 86
$ws.Cells.Item(86, 8).Value = 2108382.8  # H86: 2108398.5 -> 2108382.8
$ws.Cells.Item(86, 10).Value = 3292772.5  # J86: 3292797.2 -> 3292772.5
$ws.Cells.Item(86, 12).Value = 3292772.5  # L86: 3292797.2 -> 3292772.5
$ws.Cells.Item(86, 14).Value = -3295018.5  # N86: -3295043.2 -> -3295018.5

# Row 89
$ws.Cells.Item(89, 8).Value = 2108382.8  # H89: 2108398.5 -> 2108382.8
$ws.Cells.Item(89, 10).Value = 3292772.5  # J89: 3292797.2 -> 3292772.5
$ws.Cells.Item(89, 12).Value = 16463862.5  # L89: 16463986 -> 16463862.5
$ws.Cells.Item(89, 14).Value = -16475094.5  # N89: -16475218 -> -16475094.5

# Row 98
$ws.Cells.Item(98, 8).Value = 1213.7826  # H98: 1268.2727 -> 1213.7826
$ws.Cells.Item(98, 9).Value = 770.9  # I98: 810.6842 -> 770.9
$ws.Cells.Item(98, 11).Value = 770.9  # K98: 810.6842 -> 770.9
$ws.Cells.Item(98, 13).Value = 727.1  # M98: 687.3158 -> 727.1

# Row 99
$ws.Cells.Item(99, 8).Value = 532.55554  # H99: 614.4167 -> 532.55554
$ws.Cells.Item(99, 10).Value = 876.25  # J99: 869.2857 -> 876.25
$ws.Cells.Item(99, 12).Value = 2628.75  # L99: 2607.8571 -> 2628.75
$ws.Cells.Item(99, 14).Value = -5624.75  # N99: -5603.8571 -> -5624.75

# Row 112
$ws.Cells.Item(112, 8).Value = 2926.9092  # H112: 2824.1667 -> 2926.9092
$ws.Cells.Item(112, 10).Value = 2926.9092  # J112: 2824.1667 -> 2926.9092
$ws.Cells.Item(112, 12).Value = 8780.7276  # L112: 8472.500100000001 -> 8780.7276
$ws.Cells.Item(112, 14).Value = -10996.7276  # N112: -10688.5001 -> -10996.7276

# Row 122
$ws.Cells.Item(122, 8).Value = 1213.7826  # H122: 1268.2727 -> 1213.7826
$ws.Cells.Item(122, 9).Value = 770.9  # I122: 810.6842 -> 770.9
$ws.Cells.Item(122, 11).Value = 2312.7  # K122: 2432.0526 -> 2312.7
$ws.Cells.Item(122, 13).Value = 137.3000000000002  # M122: 17.94740000000002 -> 137.3000000000002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)

# Row 2
$ws.Cells.Item(2, 8).Value = 903.2  # H2: 853.35297 -> 903.2
$ws.Cells.Item(2, 9).Value = 909.6429000000001  # I2: 855.875 -> 909.6429000000001
$ws.Cells.Item(2, 11).Value = 909.6429000000001  # K2: 855.875 -> 909.6429000000001
$ws.Cells.Item(2, 13).Value = -796.6429000000001  # M2: -742.875 -> -796.6429000000001

# Row 25
$ws.Cells.Item(25, 8).Value = 11135  # H25: 11006 -> 11135
$ws.Cells.Item(25, 9).Value = 0  # I25: 11006 -> 0
$ws.Cells.Item(25, 10).Value = 11135  # J25: 0 -> 11135
$ws.Cells.Item(25, 11).Value = 0  # K25: 11006 -> 0
$ws.Cells.Item(25, 12).ClearContents()  # L25: 0 -> (removed)
$ws.Cells.Item(25, 13).Value = 11135  # M25: -10604 -> 11135
$ws.Cells.Item(25, 14).Value = -11939  # N25: None -> -11939

# Row 32
$ws.Cells.Item(32, 8).Value = 2082.18  # H32: 1693.11 -> 2082.18
$ws.Cells.Item(32, 9).Value = 2082.18  # I32: 1693.11 -> 2082.18
$ws.Cells.Item(32, 11).Value = 2082.18  # K32: 1693.11 -> 2082.18
$ws.Cells.Item(32, 13).Value = -1795.18  # M32: -1406.11 -> -1795.18

# Row 74
$ws.Cells.Item(74, 8).Value = 1913.8  # H74: 1774 -> 1913.8
$ws.Cells.Item(74, 9).Value = 1892.25  # I74: 1887.625 -> 1892.25
$ws.Cells.Item(74, 10).Value = 2000  # J74: 1514.2858 -> 2000
$ws.Cells.Item(74, 11).Value = 1892.25  # K74: 1887.625 -> 1892.25
$ws.Cells.Item(74, 12).Value = 2000  # L74: 1514.2858 -> 2000
$ws.Cells.Item(74, 13).Value = -1018.25  # M74: -1013.625 -> -1018.25
$ws.Cells.Item(74, 14).Value = -3748  # N74: -3262.2858 -> -3748

# Row 77
$ws.Cells.Item(77, 8).Value = 1913.8  # H77: 1774 -> 1913.8
$ws.Cells.Item(77, 9).Value = 1892.25  # I77: 1887.625 -> 1892.25
$ws.Cells.Item(77, 10).Value = 2000  # J77: 1514.2858 -> 2000
$ws.Cells.Item(77, 11).Value = 9461.25  # K77: 9438.125 -> 9461.25
$ws.Cells.Item(77, 12).Value = 10000  # L77: 7571.429 -> 10000
$ws.Cells.Item(77, 13).Value = -5093.25  # M77: -5070.125 -> -5093.25
$ws.Cells.Item(77, 14).Value = -18736  # N77: -16307.429 -> -18736

# Row 97
$ws.Cells.Item(97, 8).Value = 750.8214  # H97: 748.90625 -> 750.8214
$ws.Cells.Item(97, 9).Value = 796.53845  # I97: 809.5862 -> 796.53845
$ws.Cells.Item(97, 10).Value = 156.5  # J97: 162.33333 -> 156.5
$ws.Cells.Item(97, 11).Value = 796.53845  # K97: 809.5862 -> 796.53845
$ws.Cells.Item(97, 12).Value = 156.5  # L97: 162.33333 -> 156.5
$ws.Cells.Item(97, 13).Value = -300.53845  # M97: -313.5862 -> -300.53845
$ws.Cells.Item(97, 14).Value = -1148.5  # N97: -1154.33333 -> -1148.5

# Row 116
$ws.Cells.Item(116, 8).Value = 903.2  # H116: 853.35297 -> 903.2
$ws.Cells.Item(116, 9).Value = 909.6429000000001  # I116: 855.875 -> 909.6429000000001
$ws.Cells.Item(116, 11).Value = 909.6429000000001  # K116: 855.875 -> 909.6429000000001
$ws.Cells.Item(116, 13).Value = 1384.3571  # M116: 1438.125 -> 1384.3571

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)

# Row 3
$ws.Cells.Item(3, 8).Value = 903.2  # H3: 853.35297 -> 903.2
$ws.Cells.Item(3, 9).Value = 909.6429000000001  # I3: 855.875 -> 909.6429000000001
$ws.Cells.Item(3, 11).Value = 909.6429000000001  # K3: 855.875 -> 909.6429000000001
$ws.Cells.Item(3, 13).Value = -795.6429000000001  # M3: -741.875 -> -795.6429000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)

# Row 31
$ws.Cells.Item(31, 8).Value = 73360.86  # H31: 59486.223 -> 73360.86
$ws.Cells.Item(31, 10).Value = 337097.34  # J31: 150713.14 -> 337097.34
$ws.Cells.Item(31, 12).Value = 337097.34  # L31: 150713.14 -> 337097.34
$ws.Cells.Item(31, 14).Value = -337687.34  # N31: -151303.14 -> -337687.34

# Row 34
$ws.Cells.Item(34, 8).Value = 73360.86  # H34: 59486.223 -> 73360.86
$ws.Cells.Item(34, 10).Value = 337097.34  # J34: 150713.14 -> 337097.34
$ws.Cells.Item(34, 12).Value = 337097.34  # L34: 150713.14 -> 337097.34
$ws.Cells.Item(34, 14).Value = -337501.34  # N34: -151117.14 -> -337501.34

# Row 58
$ws.Cells.Item(58, 8).Value = 2494.4546  # H58: 2376.5833 -> 2494.4546
$ws.Cells.Item(58, 9).Value = 2271  # I58: 2151.9 -> 2271
$ws.Cells.Item(58, 11).Value = 2271  # K58: 2151.9 -> 2271
$ws.Cells.Item(58, 13).Value = -2068  # M58: -1948.9 -> -2068

# Row 94
$ws.Cells.Item(94, 8).Value = 1248.8  # H94: 1310.2222 -> 1248.8
$ws.Cells.Item(94, 10).Value = 1361  # J94: 1456 -> 1361
$ws.Cells.Item(94, 12).Value = 1361  # L94: 1456 -> 1361
$ws.Cells.Item(94, 14).Value = -2263  # N94: -2358 -> -2263

# Row 132
$ws.Cells.Item(132, 8).Value = 3252.25  # H132: 3404.4 -> 3252.25
$ws.Cells.Item(132, 9).Value = 3248.25  # I132: 3497.5 -> 3248.25
$ws.Cells.Item(132, 10).Value = 3256.25  # J132: 3342.3333 -> 3256.25
$ws.Cells.Item(132, 11).Value = 9744.75  # K132: 10492.5 -> 9744.75
$ws.Cells.Item(132, 12).Value = 9768.75  # L132: 10026.9999 -> 9768.75
$ws.Cells.Item(132, 13).Value = -7214.75  # M132: -7962.5 -> -7214.75
$ws.Cells.Item(132, 14).Value = -14828.75  # N132: -15086.9999 -> -14828.75

# Row 136
$ws.Cells.Item(136, 8).Value = 2494.4546  # H136: 2376.5833 -> 2494.4546
$ws.Cells.Item(136, 9).Value = 2271  # I136: 2151.9 -> 2271
$ws.Cells.Item(136, 11).Value = 6813  # K136: 6455.700000000001 -> 6813
$ws.Cells.Item(136, 13).Value = -4263  # M136: -3905.700000000001 -> -4263

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)

# Row 13
$ws.Cells.Item(13, 8).Value = 3618.75  # H13: 3787.5 -> 3618.75
$ws.Cells.Item(13, 9).Value = 1750  # I13: 1999.6666 -> 1750
$ws.Cells.Item(13, 10).Value = 4241.6665  # J13: 4860.2 -> 4241.6665
$ws.Cells.Item(13, 11).Value = 5250  # K13: 5998.9998 -> 5250
$ws.Cells.Item(13, 12).Value = 12724.9995  # L13: 14580.6 -> 12724.9995
$ws.Cells.Item(13, 13).Value = -5082  # M13: -5830.9998 -> -5082
$ws.Cells.Item(13, 14).Value = -13060.9995  # N13: -14916.6 -> -13060.9995

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)

# Row 2
$ws.Cells.Item(2, 8).Value = 96.40000000000001  # H2: 101.07143 -> 96.40000000000001
$ws.Cells.Item(2, 9).Value = 65.84614999999999  # I2: 68.75 -> 65.84614999999999
$ws.Cells.Item(2, 11).Value = 65.84614999999999  # K2: 68.75 -> 65.84614999999999
$ws.Cells.Item(2, 13).Value = 47.15385000000001  # M2: 44.25 -> 47.15385000000001

# Row 47
$ws.Cells.Item(47, 8).Value = 16199.8  # H47: 17142.715 -> 16199.8
$ws.Cells.Item(47, 10).Value = 16199.8  # J47: 17142.715 -> 16199.8
$ws.Cells.Item(47, 12).Value = 16199.8  # L47: 17142.715 -> 16199.8
$ws.Cells.Item(47, 14).Value = -17335.8  # N47: -18278.715 -> -17335.8

# Row 93
$ws.Cells.Item(93, 8).Value = 39960  # H93: 39959.5 -> 39960
$ws.Cells.Item(93, 10).Value = 39960  # J93: 39959.5 -> 39960
$ws.Cells.Item(93, 12).Value = 39960  # L93: 39959.5 -> 39960
$ws.Cells.Item(93, 14).Value = -43704  # N93: -43703.5 -> -43704

# Row 97
$ws.Cells.Item(97, 8).Value = 523.6667  # H97: 573.25 -> 523.6667
$ws.Cells.Item(97, 9).Value = 523.6667  # I97: 615.2 -> 523.6667
$ws.Cells.Item(97, 10).Value = 0  # J97: 447.4 -> 0
$ws.Cells.Item(97, 11).Value = 523.6667  # K97: 615.2 -> 523.6667
$ws.Cells.Item(97, 12).Value = 0  # L97: 447.4 -> 0
$ws.Cells.Item(97, 13).ClearContents()  # M97: -119.2 -> (removed)
$ws.Cells.Item(97, 14).Value = -27.66669999999999  # N97: -1439.4 -> -27.66669999999999

# Row 113
$ws.Cells.Item(113, 8).Value = 419016.12  # H113: 456853.47 -> 419016.12
$ws.Cells.Item(113, 9).Value = 626368  # I113: 668085.1 -> 626368
$ws.Cells.Item(113, 10).Value = 4312.375  # J113: 4214.143 -> 4312.375
$ws.Cells.Item(113, 11).Value = 626368  # K113: 668085.1 -> 626368
$ws.Cells.Item(113, 12).Value = 4312.375  # L113: 4214.143 -> 4312.375
$ws.Cells.Item(113, 13).Value = -624198  # M113: -665915.1 -> -624198
$ws.Cells.Item(113, 14).Value = -8652.375  # N113: -8554.143 -> -8652.375

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)

# Row 46
$ws.Cells.Item(46, 8).Value = 1868.1666  # H46: 2166.5334 -> 1868.1666
$ws.Cells.Item(46, 9).Value = 2117.0833  # I46: 2458.8 -> 2117.0833
$ws.Cells.Item(46, 10).Value = 1370.3334  # J46: 1582 -> 1370.3334
$ws.Cells.Item(46, 11).Value = 2117.0833  # K46: 2458.8 -> 2117.0833
$ws.Cells.Item(46, 12).Value = 1370.3334  # L46: 1582 -> 1370.3334
$ws.Cells.Item(46, 13).Value = -1929.0833  # M46: -2270.8 -> -1929.0833
$ws.Cells.Item(46, 14).Value = -1746.3334  # N46: -1958 -> -1746.3334

# Row 55
$ws.Cells.Item(55, 8).Value = 1500.0834  # H55: 1760.2 -> 1500.0834
$ws.Cells.Item(55, 9).Value = 161.16667  # I55: 141.6 -> 161.16667
$ws.Cells.Item(55, 10).Value = 2839  # J55: 3378.8 -> 2839
$ws.Cells.Item(55, 11).Value = 161.16667  # K55: 141.6 -> 161.16667
$ws.Cells.Item(55, 12).Value = 2839  # L55: 3378.8 -> 2839
$ws.Cells.Item(55, 13).Value = 11.83332999999999  # M55: 31.40000000000001 -> 11.83332999999999
$ws.Cells.Item(55, 14).Value = -3185  # N55: -3724.8 -> -3185

# Row 100
$ws.Cells.Item(100, 8).Value = 4749.3  # H100: 4624.8335 -> 4749.3
$ws.Cells.Item(100, 9).Value = 4599.6  # I100: 4429 -> 4599.6
$ws.Cells.Item(100, 11).Value = 4599.6  # K100: 4429 -> 4599.6
$ws.Cells.Item(100, 13).Value = -4058.6  # M100: -3888 -> -4058.6

# Row 124
$ws.Cells.Item(124, 8).Value = 99990  # H124: 99994.5 -> 99990
$ws.Cells.Item(124, 10).Value = 99990  # J124: 99994.5 -> 99990
$ws.Cells.Item(124, 12).Value = 99990  # L124: 99994.5 -> 99990
$ws.Cells.Item(124, 14).Value = -109810  # N124: -109814.5 -> -109810

# Row 132
$ws.Cells.Item(132, 8).Value = 4866  # H132: 6649.5 -> 4866

# Row 136
$ws.Cells.Item(136, 8).Value = 328553.3  # H136: 308759.16 -> 328553.3
$ws.Cells.Item(136, 10).Value = 8321.571  # J136: 7525.0625 -> 8321.571
$ws.Cells.Item(136, 12).Value = 24964.713  # L136: 22575.1875 -> 24964.713
$ws.Cells.Item(136, 14).Value = -30064.713  # N136: -27675.1875 -> -30064.713

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)

# Row 120
$ws.Cells.Item(120, 8).Value = 198200  # H120: 198300 -> 198200
$ws.Cells.Item(120, 10).Value = 198200  # J120: 198300 -> 198200
$ws.Cells.Item(120, 12).Value = 198200  # L120: 198300 -> 198200
$ws.Cells.Item(120, 14).Value = -207876  # N120: -207976 -> -207876

# Row 124
$ws.Cells.Item(124, 8).Value = 98490  # H124: 98489.664 -> 98490
$ws.Cells.Item(124, 10).Value = 98490  # J124: 98489.664 -> 98490
$ws.Cells.Item(124, 12).Value = 98490  # L124: 98489.664 -> 98490
$ws.Cells.Item(124, 14).Value = -108310  # N124: -108309.664 -> -108310

# Row 125
$ws.Cells.Item(125, 8).Value = 57932.582  # H125: 57941.168 -> 57932.582
$ws.Cells.Item(125, 10).Value = 57932.582  # J125: 57941.168 -> 57932.582
$ws.Cells.Item(125, 12).Value = 57932.582  # L125: 57941.168 -> 57932.582
$ws.Cells.Item(125, 14).Value = -67772.58199999999  # N125: -67781.16800000001 -> -67772.58199999999

# Row 132
$ws.Cells.Item(132, 8).Value = 48969.684  # H132: 47740.78 -> 48969.684
$ws.Cells.Item(132, 9).Value = 3059.3684  # I132: 3278.4119 -> 3059.3684
$ws.Cells.Item(132, 10).Value = 339735  # J132: 173717.5 -> 339735
$ws.Cells.Item(132, 11).Value = 9178.1052  # K132: 9835.235700000001 -> 9178.1052
$ws.Cells.Item(132, 12).Value = 1019205  # L132: 521152.5 -> 1019205
$ws.Cells.Item(132, 13).Value = -6648.1052  # M132: -7305.235700000001 -> -6648.1052
$ws.Cells.Item(132, 14).Value = -1024265  # N132: -526212.5 -> -1024265
